$d = $word.ActiveDocument
$d.Content.Find.Execute("46-32=", $true, $false, $false, $false, $false, $true, 1, $false, "36+39=", 1) | Out-Null
$d.Content.Find.Execute("19-9=", $true, $false, $false, $false, $false, $true, 1, $false, "52-30=", 1) | Out-Null
$d.Content.Find.Execute("62-31=", $true, $false, $false, $false, $false, $true, 1, $false, "78+16=", 1) | Out-Null
$d.Content.Find.Execute("4+76=", $true, $false, $false, $false, $false, $true, 1, $false, "50-1=", 1) | Out-Null
$d.Content.Find.Execute("46+7=", $true, $false, $false, $false, $false, $true, 1, $false, "11+18=", 1) | Out-Null
$d.Content.Find.Execute("81-22=", $true, $false, $false, $false, $false, $true, 1, $false, "2+42=", 1) | Out-Null
$d.Content.Find.Execute("98+1=", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=", 1) | Out-Null
$d.Content.Find.Execute("46+49=", $true, $false, $false, $false, $false, $true, 1, $false, "26+7=", 1) | Out-Null
$d.Content.Find.Execute("44+36=", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=", 1) | Out-Null
$d.Content.Find.Execute("27+5=", $true, $false, $false, $false, $false, $true, 1, $false, "34+44=", 1) | Out-Null
$d.Content.Find.Execute("35+4=", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=", 1) | Out-Null
$d.Content.Find.Execute("3+77=", $true, $false, $false, $false, $false, $true, 1, $false, "51-30=", 1) | Out-Null
$d.Content.Find.Execute("14+80=", $true, $false, $false, $false, $false, $true, 1, $false, "57-34=", 1) | Out-Null
$d.Content.Find.Execute("69-13=", $true, $false, $false, $false, $false, $true, 1, $false, "80-26=", 1) | Out-Null
$d.Content.Find.Execute("23-2=", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=", 1) | Out-Null
$d.Content.Find.Execute("35+55=", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=", 1) | Out-Null
$d.Content.Find.Execute("94-45=", $true, $false, $false, $false, $false, $true, 1, $false, "10+62=", 1) | Out-Null
$d.Content.Find.Execute("30+50=", $true, $false, $false, $false, $false, $true, 1, $false, "85-44=", 1) | Out-Null
$d.Content.Find.Execute("6+57=", $true, $false, $false, $false, $false, $true, 1, $false, "5+53=", 1) | Out-Null
$d.Content.Find.Execute("5+60=", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=", 1) | Out-Null
$d.Content.Find.Execute("39-6=", $true, $false, $false, $false, $false, $true, 1, $false, "65+30=", 1) | Out-Null
$d.Content.Find.Execute("32+64=", $true, $false, $false, $false, $false, $true, 1, $false, "83-75=", 1) | Out-Null
$d.Content.Find.Execute("24-2=", $true, $false, $false, $false, $false, $true, 1, $false, "50-38=", 1) | Out-Null
$d.Content.Find.Execute("76-67=", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=", 1) | Out-Null
$d.Content.Find.Execute("1+17=", $true, $false, $false, $false, $false, $true, 1, $false, "63+32=", 1) | Out-Null
$d.Content.Find.Execute("55+38=", $true, $false, $false, $false, $false, $true, 1, $false, "49-38=", 1) | Out-Null
$d.Content.Find.Execute("69-3=", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=", 1) | Out-Null
$d.Content.Find.Execute("70-51=", $true, $false, $false, $false, $false, $true, 1, $false, "6+55=", 1) | Out-Null
$d.Content.Find.Execute("9+77=", $true, $false, $false, $false, $false, $true, 1, $false, "3+53=", 1) | Out-Null
$d.Content.Find.Execute("24+25=", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=", 1) | Out-Null
$d.Content.Find.Execute("86-17=", $true, $false, $false, $false, $false, $true, 1, $false, "90-69=", 1) | Out-Null
$d.Content.Find.Execute("44+12=", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=", 1) | Out-Null
$d.Content.Find.Execute("42-32=", $true, $false, $false, $false, $false, $true, 1, $false, "41-22=", 1) | Out-Null
$d.Content.Find.Execute("32+60=", $true, $false, $false, $false, $false, $true, 1, $false, "31+56=", 1) | Out-Null
$d.Content.Find.Execute("4+71=", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=", 1) | Out-Null
$d.Content.Find.Execute("56+7=", $true, $false, $false, $false, $false, $true, 1, $false, "15-12=", 1) | Out-Null
$d.Content.Find.Execute("32+13=", $true, $false, $false, $false, $false, $true, 1, $false, "45+22=", 1) | Out-Null
$d.Content.Find.Execute("90-2=", $true, $false, $false, $false, $false, $true, 1, $false, "13+74=", 1) | Out-Null
$d.Content.Find.Execute("81-43=", $true, $false, $false, $false, $false, $true, 1, $false, "98-77=", 1) | Out-Null
$d.Content.Find.Execute("91-66=", $true, $false, $false, $false, $false, $true, 1, $false, "66+10=", 1) | Out-Null
$d.Content.Find.Execute("34+47=", $true, $false, $false, $false, $false, $true, 1, $false, "46-38=", 1) | Out-Null
$d.Content.Find.Execute("78-70=", $true, $false, $false, $false, $false, $true, 1, $false, "57-7=", 1) | Out-Null
$d.Content.Find.Execute("98-82=", $true, $false, $false, $false, $false, $true, 1, $false, "89-85=", 1) | Out-Null
$d.Content.Find.Execute("96-43=", $true, $false, $false, $false, $false, $true, 1, $false, "65-43=", 1) | Out-Null
$d.Content.Find.Execute("51-24=", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=", 1) | Out-Null
$d.Content.Find.Execute("17+22=", $true, $false, $false, $false, $false, $true, 1, $false, "53+19=", 1) | Out-Null
$d.Content.Find.Execute("28+46=", $true, $false, $false, $false, $false, $true, 1, $false, "0+72=", 1) | Out-Null
$d.Content.Find.Execute("72-57=", $true, $false, $false, $false, $false, $true, 1, $false, "31-26=", 1) | Out-Null
$d.Content.Find.Execute("21+55=", $true, $false, $false, $false, $false, $true, 1, $false, "87-13=", 1) | Out-Null
$d.Content.Find.Execute("73-67=", $true, $false, $false, $false, $false, $true, 1, $false, "8+33=", 1) | Out-Null
$d.Content.Find.Execute("46+40=", $true, $false, $false, $false, $false, $true, 1, $false, "53+43=", 1) | Out-Null
$d.Content.Find.Execute("95-78=", $true, $false, $false, $false, $false, $true, 1, $false, "6+68=", 1) | Out-Null
$d.Content.Find.Execute("7+50=", $true, $false, $false, $false, $false, $true, 1, $false, "40+56=", 1) | Out-Null
$d.Content.Find.Execute("64+5=", $true, $false, $false, $false, $false, $true, 1, $false, "60-40=", 1) | Out-Null
$d.Content.Find.Execute("37+15=", $true, $false, $false, $false, $false, $true, 1, $false, "71-25=", 1) | Out-Null
$d.Content.Find.Execute("8+65=", $true, $false, $false, $false, $false, $true, 1, $false, "37+39=", 1) | Out-Null
$d.Content.Find.Execute("68+4=", $true, $false, $false, $false, $false, $true, 1, $false, "38+43=", 1) | Out-Null
$d.Content.Find.Execute("18+61=", $true, $false, $false, $false, $false, $true, 1, $false, "80-54=", 1) | Out-Null
$d.Content.Find.Execute("85-33=", $true, $false, $false, $false, $false, $true, 1, $false, "81-55=", 1) | Out-Null
$d.Content.Find.Execute("25+65=", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=", 1) | Out-Null
$d.Content.Find.Execute("19-1=", $true, $false, $false, $false, $false, $true, 1, $false, "29+53=", 1) | Out-Null
$d.Content.Find.Execute("85-74=", $true, $false, $false, $false, $false, $true, 1, $false, "34-14=", 1) | Out-Null
$d.Content.Find.Execute("14+31=", $true, $false, $false, $false, $false, $true, 1, $false, "86-35=", 1) | Out-Null
$d.Content.Find.Execute("20+46=", $true, $false, $false, $false, $false, $true, 1, $false, "80-10=", 1) | Out-Null
$d.Content.Find.Execute("99-79=", $true, $false, $false, $false, $false, $true, 1, $false, "2+53=", 1) | Out-Null
$d.Content.Find.Execute("28+29=", $true, $false, $false, $false, $false, $true, 1, $false, "54+15=", 1) | Out-Null
$d.Content.Find.Execute("46+3=", $true, $false, $false, $false, $false, $true, 1, $false, "75-45=", 1) | Out-Null
$d.Content.Find.Execute("18+8=", $true, $false, $false, $false, $false, $true, 1, $false, "98-13=", 1) | Out-Null
$d.Content.Find.Execute("42+27=", $true, $false, $false, $false, $false, $true, 1, $false, "73-12=", 1) | Out-Null
$d.Content.Find.Execute("76-19=", $true, $false, $false, $false, $false, $true, 1, $false, "98-62=", 1) | Out-Null
$d.Content.Find.Execute("75-62=", $true, $false, $false, $false, $false, $true, 1, $false, "68+18=", 1) | Out-Null
$d.Content.Find.Execute("52-0=", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=", 1) | Out-Null
$d.Content.Find.Execute("46-44=", $true, $false, $false, $false, $false, $true, 1, $false, "54+37=", 1) | Out-Null
$d.Content.Find.Execute("65-3=", $true, $false, $false, $false, $false, $true, 1, $false, "84-33=", 1) | Out-Null
$d.Content.Find.Execute("21+64=", $true, $false, $false, $false, $false, $true, 1, $false, "11+48=", 1) | Out-Null
$d.Content.Find.Execute("65-23=", $true, $false, $false, $false, $false, $true, 1, $false, "50+42=", 1) | Out-Null
$d.Content.Find.Execute("95-13=", $true, $false, $false, $false, $false, $true, 1, $false, "70-9=", 1) | Out-Null
$d.Content.Find.Execute("57-45=", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=", 1) | Out-Null
$d.Content.Find.Execute("45+38=", $true, $false, $false, $false, $false, $true, 1, $false, "11+87=", 1) | Out-Null
$d.Content.Find.Execute("49+31=", $true, $false, $false, $false, $false, $true, 1, $false, "17+61=", 1) | Out-Null
$d.Content.Find.Execute("66-34=", $true, $false, $false, $false, $false, $true, 1, $false, "18+68=", 1) | Out-Null
$d.Content.Find.Execute("83-22=", $true, $false, $false, $false, $false, $true, 1, $false, "95-46=", 1) | Out-Null
$d.Content.Find.Execute("0+63=", $true, $false, $false, $false, $false, $true, 1, $false, "66-13=", 1) | Out-Null
$d.Content.Find.Execute("67-7=", $true, $false, $false, $false, $false, $true, 1, $false, "86-78=", 1) | Out-Null
$d.Content.Find.Execute("92-3=", $true, $false, $false, $false, $false, $true, 1, $false, "55-19=", 1) | Out-Null
$d.Content.Find.Execute("7+74=", $true, $false, $false, $false, $false, $true, 1, $false, "79-43=", 1) | Out-Null
$d.Content.Find.Execute("97-32=", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=", 1) | Out-Null
$d.Content.Find.Execute("14+68=", $true, $false, $false, $false, $false, $true, 1, $false, "51+7=", 1) | Out-Null
$d.Content.Find.Execute("45-14=", $true, $false, $false, $false, $false, $true, 1, $false, "92-46=", 1) | Out-Null
$d.Content.Find.Execute("5+13=", $true, $false, $false, $false, $false, $true, 1, $false, "66-54=", 1) | Out-Null
$d.Content.Find.Execute("55+6=", $true, $false, $false, $false, $false, $true, 1, $false, "9+45=", 1) | Out-Null
$d.Content.Find.Execute("33+25=", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=", 1) | Out-Null
$d.Content.Find.Execute("44-4=", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=", 1) | Out-Null
$d.Content.Find.Execute("30+63=", $true, $false, $false, $false, $false, $true, 1, $false, "77-43=", 1) | Out-Null
$d.Content.Find.Execute("48+9=", $true, $false, $false, $false, $false, $true, 1, $false, "13+66=", 1) | Out-Null
$d.Content.Find.Execute("92-22=", $true, $false, $false, $false, $false, $true, 1, $false, "72+12=", 1) | Out-Null
$d.Content.Find.Execute("86-33=", $true, $false, $false, $false, $false, $true, 1, $false, "51-15=", 1) | Out-Null
$d.Content.Find.Execute("54-24=", $true, $false, $false, $false, $false, $true, 1, $false, "84-6=", 1) | Out-Null
$d.Content.Find.Execute("2+49=", $true, $false, $false, $false, $false, $true, 1, $false, "17+72=", 1) | Out-Null
$d.Content.Find.Execute("0+6=", $true, $false, $false, $false, $false, $true, 1, $false, "40-14=", 1) | Out-Null
